# Update workbook per the latest metrics refresh.
$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")
$today   = $wb.Worksheets.Item("today")

# --- Metrics sheet: refreshed metric values (column B) ---
$metrics.Range("B2").Value  = 337418.4800000001
$metrics.Range("B3").Value  = 297323.29000000004
$metrics.Range("B4").Value  = 104031.5
$metrics.Range("B5").Value  = 13769
$metrics.Range("B6").Value  = 5133664.2300000014
$metrics.Range("B7").Value  = 4339399.9700000007
$metrics.Range("B8").Value  = 1510991.3300000003
$metrics.Range("B9").Value  = 199976
$metrics.Range("B10").Value = 33599045.220000014
$metrics.Range("B11").Value = 31614675.130000006
$metrics.Range("B12").Value = 11792713.369999999
$metrics.Range("B13").Value = 1297606

# --- today sheet: fill in the previously blank daily figures ---
$today.Range("B3").Value = 12181.71
$today.Range("B4").Value = 10774.28
$today.Range("B5").Value = 3728.04
$today.Range("B6").Value = 513

# --- Selections / active sheet bookkeeping ---
# "today" keeps the F11:F22 range highlighted but is no longer the tab in focus.
$today.Activate()
$today.Range("F11:F22").Select()

# "Metrics" becomes the active/selected tab with D6 selected.
$metrics.Activate()
$metrics.Range("D6").Select()
